$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New date header columns (N..U), mirroring the existing B..M pattern ---
# Merge the new header date cells first (matching the existing L1:M1-style
# pairs), THEN apply styles/values - merging re-stamps the whole range with
# the top-left cell's style, so styling afterwards keeps things correct and
# lets the engine reuse the existing cell styles (s="1" date cell, s="2"
# blank pair cell) instead of fabricating new duplicate style entries.
$ws.Range("N1:O1").Merge()
$ws.Range("P1:Q1").Merge()
$ws.Range("R1:S1").Merge()
$ws.Range("T1:U1").Merge()

$ws.Range("N1").HorizontalAlignment = $ws.Range("L1").HorizontalAlignment
$ws.Range("N1").NumberFormat = $ws.Range("L1").NumberFormat
$ws.Range("N1").Value = 41733
$ws.Range("O1").HorizontalAlignment = $ws.Range("M1").HorizontalAlignment

$ws.Range("P1").HorizontalAlignment = $ws.Range("L1").HorizontalAlignment
$ws.Range("P1").NumberFormat = $ws.Range("L1").NumberFormat
$ws.Range("P1").Value = 41734
$ws.Range("Q1").HorizontalAlignment = $ws.Range("M1").HorizontalAlignment

$ws.Range("R1").HorizontalAlignment = $ws.Range("L1").HorizontalAlignment
$ws.Range("R1").NumberFormat = $ws.Range("L1").NumberFormat
$ws.Range("R1").Value = 41735
$ws.Range("S1").HorizontalAlignment = $ws.Range("M1").HorizontalAlignment

$ws.Range("T1").HorizontalAlignment = $ws.Range("L1").HorizontalAlignment
$ws.Range("T1").NumberFormat = $ws.Range("L1").NumberFormat
$ws.Range("T1").Value = 41736
$ws.Range("U1").HorizontalAlignment = $ws.Range("M1").HorizontalAlignment

# New column widths for the newly-used columns O, P, Q, S, U
$ws.Columns.Item(15).ColumnWidth = 22.42578125 - (5/6)
$ws.Columns.Item(16).ColumnWidth = 7.5703125 - (5/6)
$ws.Columns.Item(17).ColumnWidth = 20.7109375 - (5/6)
$ws.Columns.Item(19).ColumnWidth = 16.140625 - (5/6)
$ws.Columns.Item(21).ColumnWidth = 15.5703125 - (5/6)

# New note header in row 2 for the 4/4 block
$ws.Range("O2").Value = "Leaves becoming spotty"

# New "# Leaves" data for rows 3-7
$ws.Range("N3").Value = 21
$ws.Range("P3").Value = 21
$ws.Range("R3").Value = 21
$ws.Range("T3").Value = 22

$ws.Range("N4").Value = 13
$ws.Range("P4").Value = 13
$ws.Range("R4").Value = 13
$ws.Range("T4").Value = 14

$ws.Range("N5").Value = 11
$ws.Range("P5").Value = 11
$ws.Range("Q5").Value = "Tip leaves feel dry"
$ws.Range("R5").Value = 11

$ws.Range("N6").Value = 22
$ws.Range("P6").Value = 23
$ws.Range("R6").Value = 24

$ws.Range("N7").Value = 11
$ws.Range("P7").Value = 12
$ws.Range("R7").Value = 12

# --- View state updates: scroll so column I is the left-most visible column,
# and select T5 like in the edited workbook ---
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("T5").Select()
